# RMA Receipt Reversal.xlsx - "RMA Details Maintenance Grid" sheet
# Update the RMA / Shipper line / Id test-data values for rows 2-4
# (SYDATA-Work order testcases) from the old "RMA-NSWX-*" values to the
# new "RMA-B23Z-*" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2
$ws.Range("E2").Value = "RMA-B23Z-001"
$ws.Range("F2").Value = "RMA-B23Z-1-1"
$ws.Range("J2").Value = "a7s5f000000xL3DAAU"

# Row 3
$ws.Range("E3").Value = "RMA-B23Z-002"
$ws.Range("F3").Value = "RMA-B23Z-1-2"
$ws.Range("J3").Value = "a7s5f000000xL3EAAU"

# Row 4
$ws.Range("E4").Value = "RMA-B23Z-003"
$ws.Range("F4").Value = "RMA-B23Z-1-3"
$ws.Range("J4").Value = "a7s5f000000xL3FAAU"
